$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-7 down to 7-8
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with this week's data
$ws.Range("A6").Value2 = 1
$ws.Range("B6").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value2 = "Arica y Parinacota"
$ws.Range("D6").Value2 = 45215
$ws.Range("E6").Value2 = 15
$ws.Range("F6").Value2 = 100114002
$ws.Range("G6").Value2 = "Camote"
$ws.Range("H6").Value2 = "Sin especificar"
$ws.Range("I6").Value2 = "Primera"
$ws.Range("J6").Value2 = 200
$ws.Range("K6").Value2 = 11000
$ws.Range("L6").Value2 = 12000
$ws.Range("M6").Value2 = 11500
$ws.Range("N6").Value2 = "$/malla 18 kilos"
$ws.Range("O6").Value2 = "Perú"
$ws.Range("P6").Value2 = 639
$ws.Range("Q6").Value2 = 18
$ws.Range("R6").Value2 = "Hortaliza"
